# "added form for insurance"
# Adds a second select_one question (insurance) to the survey sheet, the
# corresponding invest/insurance choice lists to the choices sheet, and
# renames the form_title in settings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 2 stays "text / email / Your email / yes" (unchanged)

# Row 3: the old "regret" yes/no question becomes the "invest" question
$survey.Range("A3").Value = "select_one invest_options"
$survey.Range("B3").Value = "invest_choice"
$survey.Range("C3").Value = "Do you want to invest?"
$survey.Range("D3").Value = "yes"

# Row 4 (new): the insurance question
$survey.Range("A4").Value = "select_one invest_options"
$survey.Range("B4").Value = "insurance_choice"
$survey.Range("C4").Value = "Do you want to buy insurance?"
$survey.Range("D4").Value = "yes"

$survey.Columns.Item(2).ColumnWidth = 32.1640625

$survey.Range("B4").Select()

# ---------------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# Rows 2-3: old yes_no list becomes the invest_options list
$choices.Range("A2").Value = "invest_options"
$choices.Range("B2").Value = "invest"
$choices.Range("C2").Value = "Invest"

$choices.Range("A3").Value = "invest_options"
$choices.Range("B3").Value = "dont_invest"
$choices.Range("C3").Value = "Don" + [char]0x2019 + "t Invest"

# Rows 4-5 (new): insurance_choice list
$choices.Range("A4").Value = "insurance_choice"
$choices.Range("B4").Value = "buy_insurance"
$choices.Range("C4").Value = "Buy Insurance"

$choices.Range("A5").Value = "insurance_choice"
$choices.Range("B5").Value = "dont_buyInsurance"
$choices.Range("C5").Value = "Dont buy Insurance"

$choices.Columns.Item(1).ColumnWidth = 25.83203125

$choices.Range("B9").Select()

# ---------------------------------------------------------------------
# settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

$settings.Range("A2").Value = "WhatsThePointGameInteractive_Section1"
# form_id / version (B2 / C2) stay "section3" / "v1"

$settings.Activate()
$settings.Range("B3").Select()

Write-Host "edit applied"
